$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.604.59"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "2.985.00"
$ws.Range("E3").Value = "  +2.74%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.95"
$ws.Range("E5").Value = "  +3.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.84"
$ws.Range("E6").Value = "  +1.49%  "

$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  +1.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.29"
$ws.Range("E10").Value = "  +2.37%  "

$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("E12").Value = "  +1.86%  "

$ws.Range("D13").Value = "3.452.23"
$ws.Range("E13").Value = "  +2.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.44"
$ws.Range("E14").Value = "  +1.08%  "

$ws.Range("E15").Value = "  +3.31%  "

$ws.Range("D16").Value = "2.981.42"
$ws.Range("E16").Value = "  +2.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.978"
$ws.Range("E17").Value = "  +5.64%  "

$ws.Range("D18").Value = "51.544.74"
$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("E19").Value = "  +3.79%  "

$ws.Range("E20").Value = "  +3.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.97"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  +2.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.39"
$ws.Range("E23").Value = "  +1.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.12"
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("E25").Value = "  +6.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.20"
$ws.Range("E26").Value = "  +17.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.72"
$ws.Range("E27").Value = "  +24.53%  "

$ws.Range("E28").Value = "  +15.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.171"
$ws.Range("E29").Value = "  -2.35%  "

$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.94"
$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.89"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.18"
$ws.Range("E33").Value = "  +2.24%  "

$ws.Range("E34").Value = "  -1.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.06"
$ws.Range("E35").Value = "  +0.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0446"
$ws.Range("E36").Value = "  +5.90%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.04"
$ws.Range("E38").Value = "  +1.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.23"
$ws.Range("E39").Value = "  +1.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.60"
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("E41").Value = "  +1.25%  "

$ws.Range("E42").Value = "  +3.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "126.08"
$ws.Range("E43").Value = "  +6.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.55"
$ws.Range("E44").Value = "  -2.55%  "

$ws.Range("E45").Value = "  +19.58%  "

$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("E47").Value = "  +3.27%  "

$ws.Range("D48").Value = "2.036.84"
$ws.Range("E48").Value = "  +1.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.25"
$ws.Range("E49").Value = "  +3.23%  "

$ws.Range("E50").Value = "  +8.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.45"
$ws.Range("E51").Value = "  +3.51%  "
